$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New record (row 6) appended to the "Artfynd" sheet.
$ws.Range("A6").Value = 112017318
$ws.Range("B6").Value = 90295
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 4740
$ws.Range("F6").Value = "Sotriska"
$ws.Range("G6").Value = "Lactarius lignyotus"
$ws.Range("H6").Value = "Fr."
$ws.Range("I6").Value = "'4"
$ws.Range("J6").Value = "fruktkroppar"
$ws.Range("P6").Value = "Granskärs våtmark, Söderhamn, Hls"
$ws.Range("Q6").Value = 610923.2678714381
$ws.Range("R6").Value = 6799713.923112066
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Gävleborg"
$ws.Range("U6").Value = "Söderhamn"
$ws.Range("V6").Value = "Hälsingland"
$ws.Range("W6").Value = "Norrala"
$ws.Range("Y6").Value = "'2023-09-10"
$ws.Range("Z6").Value = "00:00"
$ws.Range("AA6").Value = "'2023-09-10"
$ws.Range("AB6").Value = "00:00"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AI6").Value = "blandskog mossig slänt med tall och gran"
$ws.Range("AW6").Value = "Andreas Nilsson"
$ws.Range("AX6").Value = "Andreas Nilsson, Max Rosendahl"
